{"js": "// Apply the documented change:\n//  1. The placeholder empty string assigned to `url` becomes the GitHub raw\n//     CSV URL.\n//  2. The `file_name <- \"...csv\"  # remove once GitHub link is added` line is\n//     deleted (leaving a blank line in its place).\n//  3. `read_csv(file_name, ...)` becomes `read_csv(url, ...)`.\n//  4. The `sessionInfo()` \"loaded via a namespace\" package listing gets a new\n//     `curl_4.3.3` entry inserted (alphabetical reflow of the 3-per-line\n//     table, one extra wrapped line at the end).\n\nconst body = context.document.body;\n\n// --- 1. url <- \"\"  ->  url <- \"https://raw.githubusercontent.com/...\" ---\nconst urlPlaceholder = body.search('\"\"', { matchCase: true });\nurlPlaceholder.load(\"items\");\nawait context.sync();\nif (urlPlaceholder.items.length !== 1) {\n  throw new Error(\n    \"Expected exactly 1 match for the url placeholder '\\\"\\\"', found \" +\n      urlPlaceholder.items.length\n  );\n}\nurlPlaceholder.items[0].insertText(\n  '\"https://raw.githubusercontent.com/Rokkaan5/5301-project-test/main/data/PSCompPars_2022.12.07_15.29.01.csv?token=GHSAT0AAAAAAB3GS64J5Y4EH3UOGIIULXQAY4RJZ4Q\"',\n  Word.InsertLocation.replace\n);\nawait context.sync();\n\n// --- 2. Delete the whole \"file_name <- ... # remove once GitHub link is added\" line ---\nconst fileNameStart = body.search(\"file_name \", { matchCase: true });\nfileNameStart.load(\"items\");\nconst removeComment = body.search(\"# remove once GitHub link is added\", {\n  matchCase: true,\n});\nremoveComment.load(\"items\");\nawait context.sync();\nif (fileNameStart.items.length !== 1 || removeComment.items.length !== 1) {\n  throw new Error(\n    \"Expected exactly 1 match each for 'file_name ' and the trailing comment, got \" +\n      fileNameStart.items.length +\n      \"/\" +\n      removeComment.items.length\n  );\n}\nconst lineToDelete = fileNameStart.items[0].expandTo(removeComment.items[0]);\nlineToDelete.insertText(\"\", Word.InsertLocation.replace);\nawait context.sync();\n\n// --- 3. read_csv(file_name,  ->  read_csv(url, ---\nconst readCsvArg = body.search(\"(file_name,\", { matchCase: true });\nreadCsvArg.load(\"items\");\nawait context.sync();\nif (readCsvArg.items.length !== 1) {\n  throw new Error(\n    \"Expected exactly 1 match for '(file_name,', found \" + readCsvArg.items.length\n  );\n}\nreadCsvArg.items[0].insertText(\"(url,\", Word.InsertLocation.replace);\nawait context.sync();\n\n// --- 4. sessionInfo() package table: insert curl_4.3.3, reflow rest ---\nconst beforeLines = [\n  \"[13] googlesheets4_1.0.1 readxl_1.4.1        rstudioapi_0.14    \",\n  \"[16] rmarkdown_2.18      labeling_0.4.2      googledrive_2.0.0  \",\n  \"[19] bit_4.0.5           munsell_0.5.0       broom_1.0.1        \",\n  \"[22] compiler_4.2.2      modelr_0.1.10       xfun_0.35          \",\n  \"[25] pkgconfig_2.0.3     htmltools_0.5.3     tidyselect_1.2.0   \",\n  \"[28] fansi_1.0.3         crayon_1.5.2        tzdb_0.3.0         \",\n  \"[31] dbplyr_2.2.1        withr_2.5.0         grid_4.2.2         \",\n  \"[34] jsonlite_1.8.3      gtable_0.3.1        lifecycle_1.0.3    \",\n  \"[37] DBI_1.1.3           magrittr_2.0.3      scales_1.2.1       \",\n  \"[40] cli_3.4.1           stringi_1.7.8       vroom_1.6.0        \",\n  \"[43] farver_2.1.1        fs_1.5.2            xml2_1.3.3         \",\n  \"[46] ellipsis_0.3.2      generics_0.1.3      vctrs_0.5.1        \",\n  \"[49] tools_4.2.2         bit64_4.0.5         glue_1.6.2         \",\n  \"[52] hms_1.1.2           parallel_4.2.2      fastmap_1.1.0      \",\n  \"[55] yaml_2.3.6          timechange_0.1.1    colorspace_2.0-3   \",\n  \"[58] gargle_1.2.1        knitr_1.41          haven_2.5.1        \",\n];\nconst afterLines = [\n  \"[13] googlesheets4_1.0.1 curl_4.3.3          readxl_1.4.1       \",\n  \"[16] rstudioapi_0.14     rmarkdown_2.18      labeling_0.4.2     \",\n  \"[19] googledrive_2.0.0   bit_4.0.5           munsell_0.5.0      \",\n  \"[22] broom_1.0.1         compiler_4.2.2      modelr_0.1.10      \",\n  \"[25] xfun_0.35           pkgconfig_2.0.3     htmltools_0.5.3    \",\n  \"[28] tidyselect_1.2.0    fansi_1.0.3         crayon_1.5.2       \",\n  \"[31] tzdb_0.3.0          dbplyr_2.2.1        withr_2.5.0        \",\n  \"[34] grid_4.2.2          jsonlite_1.8.3      gtable_0.3.1       \",\n  \"[37] lifecycle_1.0.3     DBI_1.1.3           magrittr_2.0.3     \",\n  \"[40] scales_1.2.1        cli_3.4.1           stringi_1.7.8      \",\n  \"[43] vroom_1.6.0         farver_2.1.1        fs_1.5.2           \",\n  \"[46] xml2_1.3.3          ellipsis_0.3.2      generics_0.1.3     \",\n  \"[49] vctrs_0.5.1         tools_4.2.2         bit64_4.0.5        \",\n  \"[52] glue_1.6.2          hms_1.1.2           parallel_4.2.2     \",\n  \"[55] fastmap_1.1.0       yaml_2.3.6          timechange_0.1.1   \",\n  \"[58] colorspace_2.0-3    gargle_1.2.1        knitr_1.41         \",\n];\nconst newFinalLine = \"[61] haven_2.5.1        \";\n\n// Replace each of the 16 existing lines (in place, same run count/order)\n// with its reflowed replacement text.\nfor (let i = 0; i < beforeLines.length; i++) {\n  const lineResult = body.search(beforeLines[i], { matchCase: true });\n  lineResult.load(\"items\");\n  await context.sync();\n  if (lineResult.items.length !== 1) {\n    throw new Error(\n      \"Expected exactly 1 match for package-table line \" +\n        i +\n        \", got \" +\n        lineResult.items.length +\n        \": \" +\n        beforeLines[i]\n    );\n  }\n  lineResult.items[0].insertText(afterLines[i], Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// Append a new wrapped line (break + text, carrying the VerbatimChar style)\n// after the now-last line of the table for the newly-displaced haven_2.5.1.\nconst lastLineResult = body.search(afterLines[afterLines.length - 1], {\n  matchCase: true,\n});\nlastLineResult.load(\"items\");\nawait context.sync();\nif (lastLineResult.items.length !== 1) {\n  throw new Error(\n    \"Expected exactly 1 match for the final package-table line, got \" +\n      lastLineResult.items.length\n  );\n}\nconst lastLineRange = lastLineResult.items[0];\nlastLineRange.insertBreak(Word.BreakType.line, Word.InsertLocation.after);\nawait context.sync();\n\nconst tableParagraph = lastLineRange.paragraphs.getFirst();\nconst paragraphEnd = tableParagraph.getRange(Word.RangeLocation.end);\nconst newLineRange = paragraphEnd.insertText(\n  newFinalLine,\n  Word.InsertLocation.replace\n);\nnewLineRange.style = \"Verbatim Char\";\nawait context.sync();\n", "ps1": "# Apply the documented change:\n#  1. The placeholder empty string assigned to `url` becomes the GitHub raw\n#     CSV URL.\n#  2. The `file_name <- \"...csv\"  # remove once GitHub link is added` line is\n#     deleted (leaving a blank line in its place).\n#  3. `read_csv(file_name, ...)` becomes `read_csv(url, ...)`.\n#  4. The `sessionInfo()` \"loaded via a namespace\" package listing gets a new\n#     `curl_4.3.3` entry inserted (alphabetical reflow of the 3-per-line\n#     table, one extra wrapped line at the end).\n\n$d = $word.ActiveDocument\n\n# --- 1. url <- \"\"  ->  url <- \"https://raw.githubusercontent.com/...\" ---\n$urlRange = $d.Content\n$urlRange.Find.ClearFormatting()\n$urlRange.Find.Text = '\"\"'\n$urlRange.Find.MatchCase = $true\n$urlRange.Find.MatchWildcards = $false\n$foundUrl = $urlRange.Find.Execute()\nif (-not $foundUrl) {\n  throw \"Could not find the url placeholder '\"\"'\"\n}\n# Assign .Text directly (rather than Find.Replacement) so straight quotes\n# are not auto-corrected into curly quotes.\n$urlRange.Text = '\"https://raw.githubusercontent.com/Rokkaan5/5301-project-test/main/data/PSCompPars_2022.12.07_15.29.01.csv?token=GHSAT0AAAAAAB3GS64J5Y4EH3UOGIIULXQAY4RJZ4Q\"'\n\n# --- 2. Delete the whole \"file_name <- ... # remove once GitHub link is added\" line ---\n$startRange = $d.Content\n$startRange.Find.ClearFormatting()\n$startRange.Find.Text = \"file_name \"\n$startRange.Find.MatchCase = $true\n$foundStart = $startRange.Find.Execute()\nif (-not $foundStart) {\n  throw \"Could not find 'file_name '\"\n}\n\n$endRange = $d.Content\n$endRange.Find.ClearFormatting()\n$endRange.Find.Text = \"# remove once GitHub link is added\"\n$endRange.Find.MatchCase = $true\n$foundEnd = $endRange.Find.Execute()\nif (-not $foundEnd) {\n  throw \"Could not find '# remove once GitHub link is added'\"\n}\n\n$lineToDelete = $d.Range($startRange.Start, $endRange.End)\n$lineToDelete.Text = \"\"\n\n# --- 3. read_csv(file_name,  ->  read_csv(url, ---\n$readCsvRange = $d.Content\n$readCsvRange.Find.ClearFormatting()\n$readCsvRange.Find.Text = \"(file_name,\"\n$readCsvRange.Find.MatchCase = $true\n$foundReadCsv = $readCsvRange.Find.Execute()\nif (-not $foundReadCsv) {\n  throw \"Could not find '(file_name,'\"\n}\n$readCsvRange.Text = \"(url,\"\n\n# --- 4. sessionInfo() package table: insert curl_4.3.3, reflow rest ---\n$beforeLines = @(\n  \"[13] googlesheets4_1.0.1 readxl_1.4.1        rstudioapi_0.14    \",\n  \"[16] rmarkdown_2.18      labeling_0.4.2      googledrive_2.0.0  \",\n  \"[19] bit_4.0.5           munsell_0.5.0       broom_1.0.1        \",\n  \"[22] compiler_4.2.2      modelr_0.1.10       xfun_0.35          \",\n  \"[25] pkgconfig_2.0.3     htmltools_0.5.3     tidyselect_1.2.0   \",\n  \"[28] fansi_1.0.3         crayon_1.5.2        tzdb_0.3.0         \",\n  \"[31] dbplyr_2.2.1        withr_2.5.0         grid_4.2.2         \",\n  \"[34] jsonlite_1.8.3      gtable_0.3.1        lifecycle_1.0.3    \",\n  \"[37] DBI_1.1.3           magrittr_2.0.3      scales_1.2.1       \",\n  \"[40] cli_3.4.1           stringi_1.7.8       vroom_1.6.0        \",\n  \"[43] farver_2.1.1        fs_1.5.2            xml2_1.3.3         \",\n  \"[46] ellipsis_0.3.2      generics_0.1.3      vctrs_0.5.1        \",\n  \"[49] tools_4.2.2         bit64_4.0.5         glue_1.6.2         \",\n  \"[52] hms_1.1.2           parallel_4.2.2      fastmap_1.1.0      \",\n  \"[55] yaml_2.3.6          timechange_0.1.1    colorspace_2.0-3   \",\n  \"[58] gargle_1.2.1        knitr_1.41          haven_2.5.1        \"\n)\n$afterLines = @(\n  \"[13] googlesheets4_1.0.1 curl_4.3.3          readxl_1.4.1       \",\n  \"[16] rstudioapi_0.14     rmarkdown_2.18      labeling_0.4.2     \",\n  \"[19] googledrive_2.0.0   bit_4.0.5           munsell_0.5.0      \",\n  \"[22] broom_1.0.1         compiler_4.2.2      modelr_0.1.10      \",\n  \"[25] xfun_0.35           pkgconfig_2.0.3     htmltools_0.5.3    \",\n  \"[28] tidyselect_1.2.0    fansi_1.0.3         crayon_1.5.2       \",\n  \"[31] tzdb_0.3.0          dbplyr_2.2.1        withr_2.5.0        \",\n  \"[34] grid_4.2.2          jsonlite_1.8.3      gtable_0.3.1       \",\n  \"[37] lifecycle_1.0.3     DBI_1.1.3           magrittr_2.0.3     \",\n  \"[40] scales_1.2.1        cli_3.4.1           stringi_1.7.8      \",\n  \"[43] vroom_1.6.0         farver_2.1.1        fs_1.5.2           \",\n  \"[46] xml2_1.3.3          ellipsis_0.3.2      generics_0.1.3     \",\n  \"[49] vctrs_0.5.1         tools_4.2.2         bit64_4.0.5        \",\n  \"[52] glue_1.6.2          hms_1.1.2           parallel_4.2.2     \",\n  \"[55] fastmap_1.1.0       yaml_2.3.6          timechange_0.1.1   \",\n  \"[58] colorspace_2.0-3    gargle_1.2.1        knitr_1.41         \"\n)\n$newFinalLine = \"[61] haven_2.5.1        \"\n\nfor ($i = 0; $i -lt $beforeLines.Length; $i++) {\n  $lineRange = $d.Content\n  $lineRange.Find.ClearFormatting()\n  $lineRange.Find.Text = $beforeLines[$i]\n  $lineRange.Find.MatchCase = $true\n  $foundLine = $lineRange.Find.Execute()\n  if (-not $foundLine) {\n    throw \"Could not find package-table line $($i): $($beforeLines[$i])\"\n  }\n  $lineRange.Text = $afterLines[$i]\n}\n\n# Append a new wrapped line (break + text, carrying the VerbatimChar style)\n# after the now-last line of the table for the newly-displaced haven_2.5.1.\n$lastLineRange = $d.Content\n$lastLineRange.Find.ClearFormatting()\n$lastLineRange.Find.Text = $afterLines[$afterLines.Length - 1]\n$lastLineRange.Find.MatchCase = $true\n$foundLast = $lastLineRange.Find.Execute()\nif (-not $foundLast) {\n  throw \"Could not find the final package-table line\"\n}\n\n$lastLineRange.Collapse(0)         # wdCollapseEnd\n$lastLineRange.InsertAfter([char]11)   # manual line break (vertical tab)\n$lastLineRange.Collapse(0)\n$lastLineRange.InsertAfter($newFinalLine)\n$lastLineRange.Style = \"Verbatim Char\"\n\nWrite-Output \"done\"\n"}
